$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the date format (numFmtId 14, m/d/yyyy) already used by E2 (DOB)
# onto the three new date cells so we reuse the existing style instead of
# creating a brand-new number format.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("R2:T2").PasteSpecial(-4122) | Out-Null

# New employer / office / injury details entered into row 2
$ws.Range("K2").Value = "Test"
$ws.Range("L2").Value = 123
$ws.Range("M2").Value = "Test Address"
$ws.Range("N2").Value = "Test OfcCity"
$ws.Range("O2").Value = "Test Ofc State"
$ws.Range("P2").Value = "Test Ofc Zip"
$ws.Range("Q2").Value = 2345678909

$ws.Range("R2").Value = Get-Date -Year 2009 -Month 12 -Day 12 -Hour 0 -Minute 0 -Second 0
$ws.Range("S2").Value = Get-Date -Year 2007 -Month 12 -Day 12 -Hour 0 -Minute 0 -Second 0
$ws.Range("T2").Value = Get-Date -Year 2009 -Month 12 -Day 12 -Hour 0 -Minute 0 -Second 0

$ws.Range("U2").Value = "Test"
$ws.Range("W2").Value = "Test Illenss"
$ws.Range("X2").Value = "Eye"
$ws.Range("Y2").Value = "Yes"
$ws.Range("Z2").Value = "Yes"
$ws.Range("AA2").Value = "Yes"
$ws.Range("AB2").Value = "Yes"

# Widen the newly filled-in date columns (R:T) similarly to the other
# bestFit columns on the sheet.
$ws.Columns("R:T").ColumnWidth = 9.63

# Move the visible window / selection to the right, matching the
# author's on-screen scroll position when they finished editing.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 13
$ws.Range("AA7").Select() | Out-Null
